$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cosinor-per-day statistics (re-run of CircadiPy analysis).
# Rows 2-4 (day test=1) and rows 5-10 (day test=2) get refreshed period,
# p/q, amplitude, acrophase, mesor and zt columns. Rows 2-4 additionally
# populate p_reject/q_reject (J/K), which were blank before.

# Row 2
$ws.Range("F2").Value = 22
$ws.Range("H2").Value = 0.1440164597546433
$ws.Range("I2").Value = 0.1440164597546433
$ws.Range("J2").Value = 0.4219988657621125
$ws.Range("K2").Value = 0.4219988657621125
$ws.Range("L2").Value = 7.300777071580573
$ws.Range("M2").Value = "[-1.6452518369844, 16.246805980145545]"
$ws.Range("N2").Value = 0.1072095088975633
$ws.Range("O2").Value = 0.1072095088975633
$ws.Range("P2").Value = -1.471737098979002
$ws.Range("Q2").Value = "[-3.2076321388003883, 0.2641579408423844]"
$ws.Range("R2").Value = 0.09460028812555166
$ws.Range("S2").Value = 0.09460028812555166
$ws.Range("T2").Value = 12.13084811921328
$ws.Range("U2").Value = "[7.105696358679816, 17.155999879746744]"
$ws.Range("V2").Value = 0.00001453900044534606
$ws.Range("W2").Value = 0.00001453900044534606
$ws.Range("X2").Value = 5.153153153153152
$ws.Range("Y2").Value = -0.9249249249249232
$ws.Range("Z2").Value = 11.23123123123123

# Row 3
$ws.Range("F3").Value = 22
$ws.Range("H3").Value = 0.3551720161778551
$ws.Range("I3").Value = 0.3551720161778551
$ws.Range("J3").Value = 0.1864721989261177
$ws.Range("K3").Value = 0.1864721989261177
$ws.Range("L3").Value = 3.557274306929711
$ws.Range("M3").Value = "[-2.577868043767721, 9.692416657627144]"
$ws.Range("N3").Value = 0.2490313425720503
$ws.Range("O3").Value = 0.2490313425720503
$ws.Range("P3").Value = -1.509473947670772
$ws.Range("Q3").Value = "[-4.616474489959776, 1.5975265946182322]"
$ws.Range("R3").Value = 0.3330527417041611
$ws.Range("S3").Value = 0.3330527417041611
$ws.Range("T3").Value = 9.692078691525204
$ws.Range("U3").Value = "[6.312234550951747, 13.071922832098661]"
$ws.Range("V3").Value = 0.0000006714618774328329
$ws.Range("W3").Value = 0.0000006714618774328329
$ws.Range("X3").Value = 5.285285285285287
$ws.Range("Y3").Value = -5.593593593593591
$ws.Range("Z3").Value = 16.16416416416416

# Row 4
$ws.Range("F4").Value = 22
$ws.Range("H4").Value = 0.2967951448490087
$ws.Range("I4").Value = 0.2967951448490087
$ws.Range("J4").Value = 0.06028187891654768
$ws.Range("K4").Value = 0.06028187891654768
$ws.Range("L4").Value = 5.615779468983182
$ws.Range("M4").Value = "[-3.335005892213725, 14.566564830180088]"
$ws.Range("N4").Value = 0.2128590191419866
$ws.Range("O4").Value = 0.2128590191419866
$ws.Range("P4").Value = -0.1383684452031542
$ws.Range("Q4").Value = "[-3.0881321179431205, 2.811395227536812]"
$ws.Range("R4").Value = 0.9251487428858218
$ws.Range("S4").Value = 0.9251487428858218
$ws.Range("T4").Value = 14.58079714595538
$ws.Range("U4").Value = "[9.331791531293256, 19.829802760617497]"
$ws.Range("V4").Value = 0.000001242216467023383
$ws.Range("W4").Value = 0.000001242216467023383
$ws.Range("X4").Value = 0.4844844844844829
$ws.Range("Y4").Value = -9.843843843843846
$ws.Range("Z4").Value = 10.81281281281281

# Row 5
$ws.Range("F5").Value = 22.68000000000011
$ws.Range("H5").Value = 0.703232178777609
$ws.Range("I5").Value = 0.703232178777609
$ws.Range("L5").Value = 2.88046164099539
$ws.Range("M5").Value = "[-6.152965217010716, 11.913888499001494]"
$ws.Range("N5").Value = 0.523981698327725
$ws.Range("O5").Value = 0.523981698327725
$ws.Range("P5").Value = 0.1572368695490383
$ws.Range("Q5").Value = "[-2.9560531475219274, 3.270526886620004]"
$ws.Range("R5").Value = 0.9194289436956
$ws.Range("S5").Value = 0.9194289436956
$ws.Range("T5").Value = 13.28567502324668
$ws.Range("U5").Value = "[8.310030039289916, 18.26132000720344]"
$ws.Range("V5").Value = 0.000002589092273685978
$ws.Range("W5").Value = 0.000002589092273685978
$ws.Range("X5").Value = 22.11243243243254
$ws.Range("Y5").Value = 10.87459459459465
$ws.Range("Z5").Value = 33.35027027027043

# Row 6
$ws.Range("F6").Value = 22.68000000000011
$ws.Range("H6").Value = 0.2695561845481641
$ws.Range("I6").Value = 0.2695561845481641
$ws.Range("L6").Value = 5.212726710360004
$ws.Range("M6").Value = "[-2.9249280512897062, 13.350381472009714]"
$ws.Range("N6").Value = 0.2035816774843107
$ws.Range("O6").Value = 0.2035816774843107
$ws.Range("P6").Value = -0.08805264694746207
$ws.Range("Q6").Value = "[-3.188763714454505, 3.012658420559581]"
$ws.Range("R6").Value = 0.9546425742490989
$ws.Range("S6").Value = 0.9546425742490989
$ws.Range("T6").Value = 12.50639588314922
$ws.Range("U6").Value = "[7.870987977320247, 17.14180378897819]"
$ws.Range("V6").Value = 0.000002141725857551791
$ws.Range("W6").Value = 0.000002141725857551791
$ws.Range("X6").Value = 0.317837837837839
$ws.Range("Y6").Value = -10.87459459459465
$ws.Range("Z6").Value = 11.51027027027033

# Row 7
$ws.Range("F7").Value = 22.68000000000011
$ws.Range("H7").Value = 0.7025127118707102
$ws.Range("I7").Value = 0.7025127118707102
$ws.Range("L7").Value = 2.636493175337765
$ws.Range("M7").Value = "[-5.827940239633944, 11.100926590309474]"
$ws.Range("N7").Value = 0.5335997237445329
$ws.Range("O7").Value = 0.5335997237445329
$ws.Range("P7").Value = 0.4591316590831935
$ws.Range("Q7").Value = "[-2.6793162571156195, 3.5975795752820066]"
$ws.Range("R7").Value = 0.7696175364139082
$ws.Range("S7").Value = 0.7696175364139082
$ws.Range("T7").Value = 14.06606860005096
$ws.Range("U7").Value = "[9.562229222600347, 18.569907977501583]"
$ws.Range("V7").Value = 0.0000001156233362387127
$ws.Range("W7").Value = 0.0000001156233362387127
$ws.Range("X7").Value = 21.0227027027028
$ws.Range("Y7").Value = 9.694054054054092
$ws.Range("Z7").Value = 32.3513513513515

# Row 8
$ws.Range("F8").Value = 22.68000000000011
$ws.Range("H8").Value = 0.5850399622987384
$ws.Range("I8").Value = 0.5850399622987384
$ws.Range("L8").Value = 3.63788103895748
$ws.Range("M8").Value = "[-6.1356548216419045, 13.411416899556865]"
$ws.Range("N8").Value = 0.4573463175153076
$ws.Range("O8").Value = 0.4573463175153076
$ws.Range("P8").Value = 0.672973801669885
$ws.Range("Q8").Value = "[-2.465474114528927, 3.811421717868697]"
$ws.Range("R8").Value = 0.6678896260132721
$ws.Range("S8").Value = 0.6678896260132721
$ws.Range("T8").Value = 12.48429119084329
$ws.Range("U8").Value = "[7.50432666885631, 17.464255712830273]"
$ws.Range("V8").Value = 0.00000780862506988278
$ws.Range("W8").Value = 0.00000780862506988278
$ws.Range("X8").Value = 20.25081081081091
$ws.Range("Y8").Value = 8.922162162162207
$ws.Range("Z8").Value = 31.57945945945961

# Row 9
$ws.Range("F9").Value = 22.68000000000011
$ws.Range("H9").Value = 0.4682021383530869
$ws.Range("I9").Value = 0.4682021383530869
$ws.Range("L9").Value = 3.725575816660973
$ws.Range("M9").Value = "[-4.388259007639515, 11.83941064096146]"
$ws.Range("N9").Value = 0.360001583434701
$ws.Range("O9").Value = 0.360001583434701
$ws.Range("P9").Value = 0.748447499053424
$ws.Range("Q9").Value = "[-2.390000417145388, 3.886895415252236]"
$ws.Range("R9").Value = 0.6333286146673092
$ws.Range("S9").Value = 0.6333286146673092
$ws.Range("T9").Value = 11.49210642433685
$ws.Range("U9").Value = "[7.2366275540169624, 15.74758529465674]"
$ws.Range("V9").Value = 0.000002105139142116741
$ws.Range("W9").Value = 0.000002105139142116741
$ws.Range("X9").Value = 19.97837837837847
$ws.Range("Y9").Value = 8.649729729729772
$ws.Range("Z9").Value = 31.30702702702717

# Row 10
$ws.Range("F10").Value = 22.68000000000011
$ws.Range("H10").Value = 0.11014206753413
$ws.Range("I10").Value = 0.11014206753413
$ws.Range("L10").Value = 8.009105042274346
$ws.Range("M10").Value = "[-1.4479066222667534, 17.466116706815445]"
$ws.Range("N10").Value = 0.09495157710777491
$ws.Range("O10").Value = 0.09495157710777491
$ws.Range("P10").Value = 1.13839493553504
$ws.Range("Q10").Value = "[-0.9874475407679633, 3.2642374118380424]"
$ws.Range("R10").Value = 0.2865322495374334
$ws.Range("S10").Value = 0.2865322495374334
$ws.Range("T10").Value = 13.1560286872988
$ws.Range("U10").Value = "[7.963772933427447, 18.34828444117016]"
$ws.Range("V10").Value = 0.000006516937814948776
$ws.Range("W10").Value = 0.000006516937814948776
$ws.Range("X10").Value = 18.5708108108109
$ws.Range("Y10").Value = 10.89729729729735
$ws.Range("Z10").Value = 26.24432432432445

Write-Host "Updated cosinor_per_day statistics for rows 2-10"
